$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update group D standings (rows 13-16) with the latest round results.
# Row 13 - Santa Cruz
$ws.Range("C13").Value = 4
$ws.Range("D13").Value = 9
$ws.Range("E13").Value = 6
$ws.Range("F13").Value = 3
$ws.Range("G13").Value = 2

# Row 14 - Tropinha
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 5
$ws.Range("E14").Value = 9
$ws.Range("F14").Value = -4
$ws.Range("G14").Value = 2

# Row 15 - Am. Ma. Velho
$ws.Range("C15").Value = 4
$ws.Range("D15").Value = 8
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 5
$ws.Range("G15").Value = 2

# Row 16 - Tira Fama
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 7
$ws.Range("F16").Value = -4
$ws.Range("G16").Value = 2

# Update the active cell selection on the sheet, as reflected by the saved file.
$ws.Range("J16").Select()
